$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 90
$ws.Range("K15").Value = 0

$ws.Range("I20").Select()
$excel.ActiveWindow.ScrollRow = 7
